$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 32500
$ws.Cells.Item(3, 10).Value = 32500
$ws.Cells.Item(3, 12).Value = 32500
$ws.Cells.Item(3, 14).Value = -32728

$ws.Cells.Item(40, 8).Value = 1975.25
$ws.Cells.Item(40, 9).Value = 1975.25
$ws.Cells.Item(40, 11).Value = 1975.25
$ws.Cells.Item(40, 13).Value = -1800.25

$ws.Cells.Item(62, 8).Value = 5601.8423
$ws.Cells.Item(62, 9).Value = 3310.75
$ws.Cells.Item(62, 11).Value = 3310.75
$ws.Cells.Item(62, 13).Value = -2686.75

$ws.Cells.Item(65, 8).Value = 5601.8423
$ws.Cells.Item(65, 9).Value = 3310.75
$ws.Cells.Item(65, 11).Value = 16553.75
$ws.Cells.Item(65, 13).Value = -13433.75

$ws.Cells.Item(75, 8).Value = 49375
$ws.Cells.Item(75, 10).Value = 49375
$ws.Cells.Item(75, 12).Value = 49375
$ws.Cells.Item(75, 14).Value = -51247

$ws.Cells.Item(78, 8).Value = 49375
$ws.Cells.Item(78, 10).Value = 49375
$ws.Cells.Item(78, 12).Value = 148125
$ws.Cells.Item(78, 14).Value = -157485

$ws.Cells.Item(98, 8).Value = 706.2632
$ws.Cells.Item(98, 9).Value = 706.2632
$ws.Cells.Item(98, 11).Value = 706.2632
$ws.Cells.Item(98, 13).Value = 791.7368

$ws.Cells.Item(102, 8).Value = 32500
$ws.Cells.Item(102, 10).Value = 32500
$ws.Cells.Item(102, 12).Value = 32500
$ws.Cells.Item(102, 14).Value = -38990

$ws.Cells.Item(122, 8).Value = 706.2632
$ws.Cells.Item(122, 9).Value = 706.2632
$ws.Cells.Item(122, 11).Value = 2118.7896
$ws.Cells.Item(122, 13).Value = 331.2103999999999

$ws.Cells.Item(141, 8).Value = 4407.909
$ws.Cells.Item(141, 9).Value = 3887.6667
$ws.Cells.Item(141, 11).Value = 11663.0001
$ws.Cells.Item(141, 13).Value = -6483.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 1001
$ws.Cells.Item(22, 9).Value = 1001
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 1001
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -702
$ws.Cells.Item(22, 14).ClearContents()

$ws.Cells.Item(33, 8).Value = 15129
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 15129
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 12).Value = 15129
$ws.Cells.Item(33, 14).Value = -15787
$ws.Cells.Item(33, 13).ClearContents()

$ws.Cells.Item(41, 8).Value = 1999.25
$ws.Cells.Item(41, 9).Value = 1999.25
$ws.Cells.Item(41, 11).Value = 1999.25
$ws.Cells.Item(41, 13).Value = -1585.25

$ws.Cells.Item(44, 8).Value = 23999.4
$ws.Cells.Item(44, 10).Value = 23999.4
$ws.Cells.Item(44, 12).Value = 23999.4
$ws.Cells.Item(44, 14).Value = -24975.4

$ws.Cells.Item(97, 8).Value = 651.3570999999999
$ws.Cells.Item(97, 9).Value = 547.61536
$ws.Cells.Item(97, 11).Value = 547.61536
$ws.Cells.Item(97, 13).Value = -51.61536000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 458
$ws.Cells.Item(5, 9).Value = 367.66666
$ws.Cells.Item(5, 11).Value = 367.66666
$ws.Cells.Item(5, 13).Value = -254.66666

$ws.Cells.Item(20, 8).Value = 837791.7
$ws.Cells.Item(20, 9).Value = 7384.1665
$ws.Cells.Item(20, 10).Value = 1668199.1
$ws.Cells.Item(20, 11).Value = 7384.1665
$ws.Cells.Item(20, 12).Value = 1668199.1
$ws.Cells.Item(20, 13).Value = -7137.1665
$ws.Cells.Item(20, 14).Value = -1668693.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1089.8235
$ws.Cells.Item(16, 9).Value = 1122
$ws.Cells.Item(16, 10).Value = 1012.6
$ws.Cells.Item(16, 11).Value = 1122
$ws.Cells.Item(16, 12).Value = 1012.6
$ws.Cells.Item(16, 13).Value = -835
$ws.Cells.Item(16, 14).Value = -1586.6

$ws.Cells.Item(92, 8).Value = 25000
$ws.Cells.Item(92, 10).Value = 25000
$ws.Cells.Item(92, 12).Value = 25000
$ws.Cells.Item(92, 14).Value = -29992

$ws.Cells.Item(113, 8).Value = 1089.8235
$ws.Cells.Item(113, 9).Value = 1122
$ws.Cells.Item(113, 10).Value = 1012.6
$ws.Cells.Item(113, 11).Value = 1122
$ws.Cells.Item(113, 12).Value = 1012.6
$ws.Cells.Item(113, 13).Value = 1048
$ws.Cells.Item(113, 14).Value = -5352.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 1442.8572
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).ClearContents()

$ws.Cells.Item(34, 8).Value = 1347
$ws.Cells.Item(34, 9).Value = 853.4091
$ws.Cells.Item(34, 10).Value = 4966.6665
$ws.Cells.Item(34, 11).Value = 2560.2273
$ws.Cells.Item(34, 12).Value = 14899.9995
$ws.Cells.Item(34, 13).Value = -2476.2273
$ws.Cells.Item(34, 14).Value = -15067.9995

$ws.Cells.Item(51, 8).Value = 262.5
$ws.Cells.Item(51, 9).Value = 83.333336
$ws.Cells.Item(51, 10).Value = 800
$ws.Cells.Item(51, 11).Value = 250.000008
$ws.Cells.Item(51, 12).Value = 2400
$ws.Cells.Item(51, 13).Value = 209.999992
$ws.Cells.Item(51, 14).Value = -3320

$ws.Cells.Item(98, 8).Value = 1949.25
$ws.Cells.Item(98, 9).Value = 2374.25
$ws.Cells.Item(98, 10).Value = 1524.25
$ws.Cells.Item(98, 11).Value = 7122.75
$ws.Cells.Item(98, 12).Value = 4572.75
$ws.Cells.Item(98, 13).Value = -5624.75
$ws.Cells.Item(98, 14).Value = -7568.75

$ws.Cells.Item(122, 8).Value = 199.5
$ws.Cells.Item(122, 10).Value = 199.5
$ws.Cells.Item(122, 12).Value = 1795.5
$ws.Cells.Item(122, 14).Value = -6695.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 60000
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 60000
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 60000
$ws.Cells.Item(62, 14).Value = -61372
$ws.Cells.Item(62, 13).ClearContents()

$ws.Cells.Item(65, 8).Value = 60000
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 60000
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 180000
$ws.Cells.Item(65, 14).Value = -186864
$ws.Cells.Item(65, 13).ClearContents()

$ws.Cells.Item(70, 8).Value = 5000
$ws.Cells.Item(70, 9).Value = 5000
$ws.Cells.Item(70, 11).Value = 5000
$ws.Cells.Item(70, 13).Value = -4730

$ws.Cells.Item(73, 8).Value = 5000
$ws.Cells.Item(73, 9).Value = 5000
$ws.Cells.Item(73, 11).Value = 5000
$ws.Cells.Item(73, 13).Value = -4064

$ws.Cells.Item(100, 8).Value = 45624.6
$ws.Cells.Item(100, 10).Value = 45624.6
$ws.Cells.Item(100, 12).Value = 45624.6
$ws.Cells.Item(100, 14).Value = -47788.6

$ws.Cells.Item(102, 8).Value = 3039.9285
$ws.Cells.Item(102, 9).Value = 1008.4286
$ws.Cells.Item(102, 10).Value = 5071.4287
$ws.Cells.Item(102, 11).Value = 1008.4286
$ws.Cells.Item(102, 12).Value = 5071.4287
$ws.Cells.Item(102, 13).Value = 613.5714
$ws.Cells.Item(102, 14).Value = -8315.4287

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 10499
$ws.Cells.Item(93, 9).Value = 10499
$ws.Cells.Item(93, 11).Value = 10499
$ws.Cells.Item(93, 13).Value = -9251

$ws.Cells.Item(122, 8).Value = 2459.8333
$ws.Cells.Item(122, 9).Value = 2101.8
$ws.Cells.Item(122, 10).Value = 4250
$ws.Cells.Item(122, 11).Value = 6305.400000000001
$ws.Cells.Item(122, 12).Value = 12750
$ws.Cells.Item(122, 13).Value = -3855.400000000001
$ws.Cells.Item(122, 14).Value = -17650

$ws.Cells.Item(132, 8).Value = 348023.66
$ws.Cells.Item(132, 9).Value = 502035.5
$ws.Cells.Item(132, 10).Value = 40000
$ws.Cells.Item(132, 11).Value = 1506106.5
$ws.Cells.Item(132, 12).Value = 120000
$ws.Cells.Item(132, 13).Value = -1503576.5
$ws.Cells.Item(132, 14).Value = -125060

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 49999
$ws.Cells.Item(74, 10).Value = 49999
$ws.Cells.Item(74, 12).Value = 49999
$ws.Cells.Item(74, 14).Value = -51871

$ws.Cells.Item(77, 8).Value = 49999
$ws.Cells.Item(77, 10).Value = 49999
$ws.Cells.Item(77, 12).Value = 149997
$ws.Cells.Item(77, 14).Value = -159357

$ws.Cells.Item(96, 8).Value = 4503
$ws.Cells.Item(96, 9).Value = 4503
$ws.Cells.Item(96, 11).Value = 4503
$ws.Cells.Item(96, 13).Value = -3130

$ws.Cells.Item(100, 8).Value = 1985.7142
$ws.Cells.Item(100, 9).Value = 1739.4
$ws.Cells.Item(100, 11).Value = 3478.8
$ws.Cells.Item(100, 13).Value = -2937.8

$ws.Cells.Item(132, 8).Value = 2504
$ws.Cells.Item(132, 9).Value = 2504
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 7512
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -4982
$ws.Cells.Item(132, 14).ClearContents()
